$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (C) column date value for all data rows (2-357) from 45184 to 45186
for ($r = 2; $r -le 357; $r++) {
    $ws.Cells.Item($r, 3).Value = 45186
}

# Add a friendly display text (the designation) as the second argument to existing HYPERLINK formulas

# Row 2: A 14414-2020
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/artfynd/A 14414-2020.xlsx", "A 14414-2020")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/kartor/A 14414-2020.png", "A 14414-2020")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/knärot/A 14414-2020.png", "A 14414-2020")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomål/A 14414-2020.docx", "A 14414-2020")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomålsmail/A 14414-2020.docx", "A 14414-2020")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsyn/A 14414-2020.docx", "A 14414-2020")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsynsmail/A 14414-2020.docx", "A 14414-2020")'

# Row 3: A 11376-2019
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/artfynd/A 11376-2019.xlsx", "A 11376-2019")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/kartor/A 11376-2019.png", "A 11376-2019")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomål/A 11376-2019.docx", "A 11376-2019")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomålsmail/A 11376-2019.docx", "A 11376-2019")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsyn/A 11376-2019.docx", "A 11376-2019")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsynsmail/A 11376-2019.docx", "A 11376-2019")'

# Row 4: A 54706-2021
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/artfynd/A 54706-2021.xlsx", "A 54706-2021")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/kartor/A 54706-2021.png", "A 54706-2021")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomål/A 54706-2021.docx", "A 54706-2021")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomålsmail/A 54706-2021.docx", "A 54706-2021")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsyn/A 54706-2021.docx", "A 54706-2021")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsynsmail/A 54706-2021.docx", "A 54706-2021")'

# Row 5: A 57791-2021
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/artfynd/A 57791-2021.xlsx", "A 57791-2021")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/kartor/A 57791-2021.png", "A 57791-2021")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomål/A 57791-2021.docx", "A 57791-2021")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomålsmail/A 57791-2021.docx", "A 57791-2021")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsyn/A 57791-2021.docx", "A 57791-2021")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsynsmail/A 57791-2021.docx", "A 57791-2021")'

# Row 6: A 67507-2021
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/artfynd/A 67507-2021.xlsx", "A 67507-2021")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/kartor/A 67507-2021.png", "A 67507-2021")'
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomål/A 67507-2021.docx", "A 67507-2021")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomålsmail/A 67507-2021.docx", "A 67507-2021")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsyn/A 67507-2021.docx", "A 67507-2021")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsynsmail/A 67507-2021.docx", "A 67507-2021")'

# Row 7: A 18498-2019
$ws.Range("S7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/artfynd/A 18498-2019.xlsx", "A 18498-2019")'
$ws.Range("T7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/kartor/A 18498-2019.png", "A 18498-2019")'
$ws.Range("V7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomål/A 18498-2019.docx", "A 18498-2019")'
$ws.Range("W7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomålsmail/A 18498-2019.docx", "A 18498-2019")'
$ws.Range("X7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsyn/A 18498-2019.docx", "A 18498-2019")'
$ws.Range("Y7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsynsmail/A 18498-2019.docx", "A 18498-2019")'

# Row 8: A 37323-2020
$ws.Range("S8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/artfynd/A 37323-2020.xlsx", "A 37323-2020")'
$ws.Range("T8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/kartor/A 37323-2020.png", "A 37323-2020")'
$ws.Range("V8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomål/A 37323-2020.docx", "A 37323-2020")'
$ws.Range("W8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomålsmail/A 37323-2020.docx", "A 37323-2020")'
$ws.Range("X8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsyn/A 37323-2020.docx", "A 37323-2020")'
$ws.Range("Y8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsynsmail/A 37323-2020.docx", "A 37323-2020")'

# Row 9: A 50811-2020
$ws.Range("S9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/artfynd/A 50811-2020.xlsx", "A 50811-2020")'
$ws.Range("T9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/kartor/A 50811-2020.png", "A 50811-2020")'
$ws.Range("V9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomål/A 50811-2020.docx", "A 50811-2020")'
$ws.Range("W9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomålsmail/A 50811-2020.docx", "A 50811-2020")'
$ws.Range("X9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsyn/A 50811-2020.docx", "A 50811-2020")'
$ws.Range("Y9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsynsmail/A 50811-2020.docx", "A 50811-2020")'

# Row 10: A 56400-2020
$ws.Range("S10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/artfynd/A 56400-2020.xlsx", "A 56400-2020")'
$ws.Range("T10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/kartor/A 56400-2020.png", "A 56400-2020")'
$ws.Range("V10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomål/A 56400-2020.docx", "A 56400-2020")'
$ws.Range("W10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomålsmail/A 56400-2020.docx", "A 56400-2020")'
$ws.Range("X10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsyn/A 56400-2020.docx", "A 56400-2020")'
$ws.Range("Y10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsynsmail/A 56400-2020.docx", "A 56400-2020")'

# Row 11: A 64046-2020
$ws.Range("S11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/artfynd/A 64046-2020.xlsx", "A 64046-2020")'
$ws.Range("T11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/kartor/A 64046-2020.png", "A 64046-2020")'
$ws.Range("V11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomål/A 64046-2020.docx", "A 64046-2020")'
$ws.Range("W11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomålsmail/A 64046-2020.docx", "A 64046-2020")'
$ws.Range("X11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsyn/A 64046-2020.docx", "A 64046-2020")'
$ws.Range("Y11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsynsmail/A 64046-2020.docx", "A 64046-2020")'

# Row 12: A 13606-2022
$ws.Range("S12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/artfynd/A 13606-2022.xlsx", "A 13606-2022")'
$ws.Range("T12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/kartor/A 13606-2022.png", "A 13606-2022")'
$ws.Range("V12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomål/A 13606-2022.docx", "A 13606-2022")'
$ws.Range("W12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomålsmail/A 13606-2022.docx", "A 13606-2022")'
$ws.Range("X12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsyn/A 13606-2022.docx", "A 13606-2022")'
$ws.Range("Y12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsynsmail/A 13606-2022.docx", "A 13606-2022")'

# Row 13: A 8144-2023
$ws.Range("S13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/artfynd/A 8144-2023.xlsx", "A 8144-2023")'
$ws.Range("T13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/kartor/A 8144-2023.png", "A 8144-2023")'
$ws.Range("V13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomål/A 8144-2023.docx", "A 8144-2023")'
$ws.Range("W13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomålsmail/A 8144-2023.docx", "A 8144-2023")'
$ws.Range("X13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsyn/A 8144-2023.docx", "A 8144-2023")'
$ws.Range("Y13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsynsmail/A 8144-2023.docx", "A 8144-2023")'

# Row 95: A 43341-2019
$ws.Range("U95").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/knärot/A 43341-2019.png", "A 43341-2019")'
$ws.Range("V95").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomål/A 43341-2019.docx", "A 43341-2019")'
$ws.Range("W95").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/klagomålsmail/A 43341-2019.docx", "A 43341-2019")'
$ws.Range("X95").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsyn/A 43341-2019.docx", "A 43341-2019")'
$ws.Range("Y95").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MARIESTAD/tillsynsmail/A 43341-2019.docx", "A 43341-2019")'
